# Add a new "Save" column (column H) to the s_vals sheet, matching the
# formatting of the existing header cells and filling data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (G1) onto
# the new header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill the new "Save" data column with 0 for each existing data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
